# Refresh the crypto price/time snapshot (column D = Price, column G = Hora)
# for data rows 2-51, as pulled from the latest coinranking.com scrape.
#
# Both columns are stored as text in the sheet (not numbers), so each new
# value is entered with a leading apostrophe to force text entry and avoid
# Excel auto-converting these numeric-looking strings into real numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Price" (column D) values, keyed by row number. Rows not listed here
# keep their existing Price value unchanged.
$priceUpdates = @{
    2  = "243.27"
    4  = "5.270"
    5  = "0.05822"
    6  = "6.461"
    7  = "3.336"
    8  = "0.8086"
    9  = "0.9025"
    10 = "0.1376"
    11 = "0.07101"
    12 = "0.03085"
    13 = "0.03028"
    14 = "0.09319"
    15 = "3.826"
    16 = "0.001537"
    17 = "0.04710"
    18 = "0.0006011"
    19 = "0.006231"
    21 = "0.003883"
    22 = "0.00008705"
    23 = "3.561"
    24 = "2.170"
    25 = "0.3191"
    40 = "0.03789"
    41 = "0.006294"
    42 = "0.1052"
    43 = "0.002503"
    44 = "0.006937"
    45 = "0.00005313"
    46 = "0.00000000750"
    47 = "0.5111"
    48 = "0.007338"
    49 = "0.00002100"
    50 = "0.0002000"
}

foreach ($row in $priceUpdates.Keys) {
    $ws.Range("D$row").Value = "'" + $priceUpdates[$row]
}

# "Hora" (column G) moves from 16 to 17 for every data row (2-51).
for ($row = 2; $row -le 51; $row++) {
    $ws.Range("G$row").Value = "'17"
}
